$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the consumption/generation/invoice/credit/balance headers to their
# shorter forms (dropping the unit suffixes), and re-order the "Credito"/
# "Fatura" columns' underlying shared-string insertion order so G1 ("Credito")
# is registered after H1 ("Fatura").
$ws.Range("E1").Value = "Consumo"
$ws.Range("F1").Value = "Geração"
$ws.Range("H1").Value = "Fatura"
$ws.Range("G1").Value = "Credito"
$ws.Range("I1").Value = "Saldo"

# Move the active selection to I1 (was A3).
$ws.Range("I1").Select() | Out-Null
